$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: mark Invalid and Absent
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: mark Total Attendance Count and Real
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: mark Total Attendance Count and Real
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: mark Absent
$ws.Range("H6").Value = 1

# Row 7: mark Absent
$ws.Range("H7").Value = 1

# Row 8: mark Absent
$ws.Range("H8").Value = 1

# Row 9: mark Absent
$ws.Range("H9").Value = 1

# Row 10: mark Total Attendance Count and Real
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

# Row 11: mark Absent
$ws.Range("H11").Value = 1

# Row 12: mark Absent
$ws.Range("H12").Value = 1

# Row 13: mark Absent
$ws.Range("H13").Value = 1

# Row 14: mark Absent
$ws.Range("H14").Value = 1

# Row 15: mark Absent
$ws.Range("H15").Value = 1

# Row 16: mark Absent
$ws.Range("H16").Value = 1

# Row 17: mark Absent
$ws.Range("H17").Value = 1

# Row 18: mark Absent
$ws.Range("H18").Value = 1
